# Localization: Content - All languages enabled! (Except Android Chinese :( )
#
# Enable the "android" (F) and "iOS" (G) columns for every localization
# entry on the "tech" sheet, except keep Android disabled for simplified
# Chinese (row 12 / sku "lang_chinese").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tech")

# Rows 6-16 (French, Italian, German, Spanish, Brazilian, Russian,
# Chinese (simplified), Japanese, Korean, Chinese (traditional), Turkish)
# Set both android (F) and iOS (G) to TRUE...
for ($r = 6; $r -le 16; $r++) {
    $ws.Cells.Item($r, 6).Value = $true
    $ws.Cells.Item($r, 7).Value = $true
}

# ...except row 12 (lang_chinese / simplified Chinese) keeps Android (F) off.
$ws.Cells.Item(12, 6).Value = $false
